$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4975.1665
$ws.Range("I74").Value = 4968.6665
$ws.Range("J74").Value = 4977.3335
$ws.Range("K74").Value = 4968.6665
$ws.Range("L74").Value = 4977.3335
$ws.Range("M74").Value = -4032.6665
$ws.Range("N74").Value = -6849.3335
$ws.Range("H77").Value = 4975.1665
$ws.Range("I77").Value = 4968.6665
$ws.Range("J77").Value = 4977.3335
$ws.Range("K77").Value = 24843.3325
$ws.Range("L77").Value = 24886.6675
$ws.Range("M77").Value = -20163.3325
$ws.Range("N77").Value = -34246.6675
$ws.Range("H98").Value = 3717.2632
$ws.Range("I98").Value = 2572.258
$ws.Range("J98").Value = 8788
$ws.Range("K98").Value = 2572.258
$ws.Range("L98").Value = 8788
$ws.Range("M98").Value = -1074.258
$ws.Range("N98").Value = -11784
$ws.Range("H113").Value = 3339
$ws.Range("I113").Value = 2505
$ws.Range("J113").Value = 3895
$ws.Range("K113").Value = 2505
$ws.Range("L113").Value = 3895
$ws.Range("M113").Value = 749
$ws.Range("N113").Value = -10403
$ws.Range("H122").Value = 3717.2632
$ws.Range("I122").Value = 2572.258
$ws.Range("J122").Value = 8788
$ws.Range("K122").Value = 7716.773999999999
$ws.Range("L122").Value = 26364
$ws.Range("M122").Value = -5266.773999999999
$ws.Range("N122").Value = -31264
$ws.Range("H129").Value = 1067.9219
$ws.Range("I129").Value = 403.72726
$ws.Range("J129").Value = 1205.7736
$ws.Range("K129").Value = 1211.18178
$ws.Range("L129").Value = 3617.3208
$ws.Range("M129").Value = 3788.81822
$ws.Range("N129").Value = -13617.3208
$ws.Range("H135").Value = 781.0857
$ws.Range("I135").Value = 723.871
$ws.Range("J135").Value = 1224.5
$ws.Range("K135").Value = 6514.839
$ws.Range("L135").Value = 11020.5
$ws.Range("M135").Value = -3979.839
$ws.Range("N135").Value = -16090.5
$ws.Range("H138").Value = 2307.2898
$ws.Range("I138").Value = 1288.3914
$ws.Range("K138").Value = 3865.1742
$ws.Range("M138").Value = 1274.8258
$ws.Range("H141").Value = 5169.8066
$ws.Range("I141").Value = 1836.2142
$ws.Range("J141").Value = 36283.332
$ws.Range("K141").Value = 5508.642599999999
$ws.Range("L141").Value = 108849.996
$ws.Range("M141").Value = -328.6425999999992
$ws.Range("N141").Value = -119209.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 13100
$ws.Range("I3").Value = 16180
$ws.Range("J3").Value = 9250
$ws.Range("K3").Value = 16180
$ws.Range("L3").Value = 9250
$ws.Range("M3").Value = -16065
$ws.Range("N3").Value = -9480
$ws.Range("H74").Value = 1245.9333
$ws.Range("I74").Value = 1117.5
$ws.Range("J74").Value = 1502.8
$ws.Range("K74").Value = 1117.5
$ws.Range("L74").Value = 1502.8
$ws.Range("M74").Value = -243.5
$ws.Range("N74").Value = -3250.8
$ws.Range("H77").Value = 1245.9333
$ws.Range("I77").Value = 1117.5
$ws.Range("J77").Value = 1502.8
$ws.Range("K77").Value = 5587.5
$ws.Range("L77").Value = 7514
$ws.Range("M77").Value = -1219.5
$ws.Range("N77").Value = -16250
$ws.Range("H97").Value = 1051.75
$ws.Range("I97").Value = 716.25
$ws.Range("K97").Value = 716.25
$ws.Range("M97").Value = -220.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 202769.2
$ws.Range("J86").Value = 401597.4
$ws.Range("L86").Value = 401597.4
$ws.Range("N86").Value = -403843.4
$ws.Range("H89").Value = 202769.2
$ws.Range("J89").Value = 401597.4
$ws.Range("L89").Value = 2007987
$ws.Range("N89").Value = -2019219
$ws.Range("H94").Value = 53927.473
$ws.Range("I94").Value = 604.3
$ws.Range("K94").Value = 604.3
$ws.Range("M94").Value = -153.3
$ws.Range("H134").Value = 2952.3
$ws.Range("I134").Value = 2652.875
$ws.Range("K134").Value = 7958.625
$ws.Range("M134").Value = -5423.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 288.58823
$ws.Range("J22").Value = 507.33334
$ws.Range("L22").Value = 507.33334
$ws.Range("N22").Value = -1207.33334
$ws.Range("H31").Value = 2198.0571
$ws.Range("I31").Value = 1331.9565
$ws.Range("K31").Value = 1331.9565
$ws.Range("M31").Value = -1036.9565
$ws.Range("H34").Value = 2198.0571
$ws.Range("I34").Value = 1331.9565
$ws.Range("K34").Value = 1331.9565
$ws.Range("M34").Value = -1129.9565
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H88").Value = 35476.57
$ws.Range("J88").Value = 35476.57
$ws.Range("L88").Value = 35476.57
$ws.Range("N88").Value = -36288.57
$ws.Range("H91").Value = 35476.57
$ws.Range("J91").Value = 35476.57
$ws.Range("L91").Value = 35476.57
$ws.Range("N91").Value = -38284.57
$ws.Range("H118").Value = 25000
$ws.Range("J118").Value = 25000
$ws.Range("L118").Value = 25000
$ws.Range("N118").Value = -28314
$ws.Range("H132").Value = 302033.12
$ws.Range("I132").Value = 398561.6
$ws.Range("J132").Value = 3672.3635
$ws.Range("K132").Value = 1195684.8
$ws.Range("L132").Value = 11017.0905
$ws.Range("M132").Value = -1193154.8
$ws.Range("N132").Value = -16077.0905
$ws.Range("H134").Value = 1615.7715
$ws.Range("I134").Value = 1070.36
$ws.Range("J134").Value = 2979.3
$ws.Range("K134").Value = 3211.08
$ws.Range("L134").Value = 8937.900000000001
$ws.Range("M134").Value = -676.0799999999999
$ws.Range("N134").Value = -14007.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 42260.082
$ws.Range("I121").Value = 398.66666
$ws.Range("J121").Value = 167844.33
$ws.Range("K121").Value = 1195.99998
$ws.Range("L121").Value = 503532.99
$ws.Range("M121").Value = 114.0000199999999
$ws.Range("N121").Value = -506152.99
$ws.Range("H131").Value = 2883.0789
$ws.Range("I131").Value = 12901.25
$ws.Range("J131").Value = 1704.4706
$ws.Range("K131").Value = 38703.75
$ws.Range("L131").Value = 5113.4118
$ws.Range("M131").Value = -33663.75
$ws.Range("N131").Value = -15193.4118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 214890.14
$ws.Range("I97").Value = 125615
$ws.Range("J97").Value = 333923.66
$ws.Range("K97").Value = 125615
$ws.Range("L97").Value = 333923.66
$ws.Range("M97").Value = -125119
$ws.Range("N97").Value = -334915.66
$ws.Range("H107").Value = 1008.2105
$ws.Range("I107").Value = 1043.8462
$ws.Range("K107").Value = 1043.8462
$ws.Range("M107").Value = 876.1538
$ws.Range("H132").Value = 2398.24
$ws.Range("I132").Value = 2064.3096
$ws.Range("J132").Value = 4151.375
$ws.Range("K132").Value = 6192.9288
$ws.Range("L132").Value = 12454.125
$ws.Range("M132").Value = -3662.9288
$ws.Range("N132").Value = -17514.125
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -39920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 21836.908
$ws.Range("I61").Value = 26067.334
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 26067.334
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -25865.334
$ws.Range("N61").Value = -3204
$ws.Range("H113").Value = 21836.908
$ws.Range("I113").Value = 26067.334
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 26067.334
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = -23897.334
$ws.Range("N113").Value = -7140
$ws.Range("H132").Value = 3100.139
$ws.Range("I132").Value = 2408.28
$ws.Range("J132").Value = 4672.5454
$ws.Range("K132").Value = 7224.84
$ws.Range("L132").Value = 14017.6362
$ws.Range("M132").Value = -4694.84
$ws.Range("N132").Value = -19077.6362
$ws.Range("H136").Value = 19805762
$ws.Range("I136").Value = 28572868
$ws.Range("J136").Value = 627717.9399999999
$ws.Range("K136").Value = 85718604
$ws.Range("L136").Value = 1883153.82
$ws.Range("M136").Value = -85716054
$ws.Range("N136").Value = -1888253.82

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1941.7333
$ws.Range("I96").Value = 1375.1428
$ws.Range("K96").Value = 1375.1428
$ws.Range("M96").Value = -2.142800000000079
$ws.Range("H132").Value = 1743.6666
$ws.Range("I132").Value = 1254.2916
$ws.Range("K132").Value = 3762.8748
$ws.Range("M132").Value = -1232.8748
